# Add a new "2021" column (R) to the data table, mirroring the formatting
# of the existing "2020" column (Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column Q (rows 3-14) into column R so the new
# column matches the existing table styling (borders, number format, etc.)
$ws.Range("Q3:Q14").Copy()
$ws.Range("R3:R14").PasteSpecial(-4122) # xlPasteFormats

# New header value for 2021
$ws.Range("R4").Value = 2021

# New data values for 2021
$ws.Range("R5").Value = 33.299999999999997
$ws.Range("R6").Value = 38.299999999999997
$ws.Range("R7").Value = 31.7
$ws.Range("R8").Value = 98.7
$ws.Range("R9").Value = 157.19999999999999
$ws.Range("R10").Value = 24.9
$ws.Range("R11").Value = 38.4
$ws.Range("R12").Value = 15.1
$ws.Range("R13").Value = 14.6
$ws.Range("R14").Value = 21.7

# Update the view selection state to mirror what a user would see after
# this edit (clicked on S6 last).
$ws.Range("S6").Select()
